$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Hung's (row 4 / column B) attendance note is a running, appended log.
# The previous last line "04/04: Ngu quyen buoi chieu" is split so that
# "chieu" starts a brand-new dated entry for 07/04.
$newNote = "10/03: Xin về sớm`n13/03: Đi trễ`n16/03: Đi trễ`n16/03: Xin về sơm`n23/03: Chưa chuẩn bị bài`n29/03: Chưa làm bài`n30/03: Xin về sớm`n04/04: Ngủ quyên buổi`n07/04: đi trể chiều"
$ws.Range("B4").Value = $newNote

# Leave the view the way the author left it: scrolled down so row 7 is at
# the top, with the active selection on G17.
$win = $excel.ActiveWindow
$win.ScrollRow = 7
$win.ScrollColumn = 1
$ws.Range("G17").Select()

$wb.Save()
